function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

function Replace-ParagraphXml($doc, $paraIndex, $innerXml) {
    $rng = $doc.Paragraphs($paraIndex).Range
    $pkg = '<?xml version="1.0" standalone="yes"?>' + `
           '<?mso-application progid="Word.Document"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body>' + $innerXml + '</w:body>' + `
           '</w:document>' + `
           '</pkg:xmlData>' + `
           '</pkg:part>' + `
           '</pkg:package>'
    $rng.InsertXML($pkg)
}

$d = $word.ActiveDocument

# --- Change 1: "Name: register" heading becomes "Name: " / "User " / "Managent"
#     (a deliberate misspelling, flagged with proofErr spellcheck markers) ---
$idx1 = Find-ParagraphIndex $d "Name: register"
$p1 = '<w:p>' + `
        '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">Name: </w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">User </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>Managent</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
      '</w:p>'
Replace-ParagraphXml $d $idx1 $p1

# --- Change 2: drop the standalone "Name: Authentication" heading paragraph
#     that preceded "url: /login" ---
$idx2 = Find-ParagraphIndex $d "Name: Authentication"
$d.Paragraphs($idx2).Range.Delete()

# --- Change 3: split the "Content: { " run into "Content" / ": { ", adding a
#     _GoBack bookmark that spans from the split point to the end of the paragraph ---
$idx3 = Find-ParagraphIndex $d "JWT token}"
$p3 = '<w:p>' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>' + `
        '<w:r><w:t>Code: 200 </w:t></w:r>' + `
        '<w:r><w:br/><w:t>Content</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:r><w:t xml:space="preserve">: { </w:t></w:r>' + `
        '<w:r><w:t>Authentication</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> : </w:t></w:r>' + `
        '<w:r><w:t>JWT token}</w:t></w:r>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
      '</w:p>'
Replace-ParagraphXml $d $idx3 $p3

# --- Change 4: remove the now-superfluous _GoBack bookmark pair that trails the
#     bare "{" run (the bookmark moved to change 3's paragraph). Identify the
#     paragraph by content rather than index, since change 2 shifted indices;
#     strip the trailing paragraph mark (\r) before testing EndsWith ---
$idx4 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($t.Contains("Code:") -and $t.Contains("Content:") -and $t.EndsWith("{")) {
        $idx4 = $i
    }
}
$p4 = '<w:p>' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>' + `
        '<w:r><w:t>Code:</w:t></w:r>' + `
        '<w:r><w:t> 200</w:t></w:r>' + `
        '<w:r><w:br/><w:t>Content: </w:t></w:r>' + `
        '<w:r><w:t>{</w:t></w:r>' + `
      '</w:p>'
Replace-ParagraphXml $d $idx4 $p4
